# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (A1) to the new headers
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-47 get the same team record values
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 92  # AD
    $ws.Cells.Item($r, 31).Value = 70  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
